{"js": "// Diary entry for \"6. Mai\": the single-word paragraph \"Tutorial\" is expanded\n// into a full diary entry, authored as three separate text runs (matching\n// how the original author's edit was captured in the source OOXML).\nconst run1 = \"Das Tutorial ist fertig. Am Anfang gab es ein paar Probleme um ein Tutorial Objekt zu erstellen, da es schwierig zu kontrollieren ist, wann genau Swing Objekte erstellt\";\nconst run2 = \" und die Tutorial Klasse danach f\u00fcr paar Sekunden schlafen sollte. Dieses Problem wurde mit dem Timer gel\u00f6st. Der Timer stellte sich im Nachhinein auch in anderen Klassen als sehr n\u00fctzlich heraus um Tasks sp\u00e4ter auszuf\u00fchren. \";\nconst run3 = \"Videos k\u00f6nnen leider immernoch nicht abgespielt werden.\";\n\n// Find the paragraph whose whole text is exactly \"Tutorial\" (the short\n// diary placeholder written on \"6. Mai\").\nconst results = context.document.body.search(\"Tutorial\", { matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Tutorial\" text to replace.');\n}\n\nconst target = results.items[0];\n\n// Build a tiny OOXML fragment with three independent <w:r> runs (no run\n// properties, same as the original single run) and insert it right before\n// the existing \"Tutorial\" text. Using insertOoxml keeps the three runs\n// distinct instead of Word's usual same-formatting run-merging.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n  '<w:r><w:t>' + run1 + '</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">' + run2 + '</w:t></w:r>' +\n  '<w:r><w:t>' + run3 + '</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n\n// Now clear out the original \"Tutorial\" text, leaving the bookmark\n// (_GoBack) that immediately follows it untouched and in place.\ntarget.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Diary entry for \"6. Mai\": the single-word paragraph \"Tutorial\" is expanded\n# into a full diary entry, authored as three separate text runs (matching\n# how the original author's edit was captured in the source OOXML).\n$run1 = 'Das Tutorial ist fertig. Am Anfang gab es ein paar Probleme um ein Tutorial Objekt zu erstellen, da es schwierig zu kontrollieren ist, wann genau Swing Objekte erstellt'\n$run2 = ' und die Tutorial Klasse danach f\u00fcr paar Sekunden schlafen sollte. Dieses Problem wurde mit dem Timer gel\u00f6st. Der Timer stellte sich im Nachhinein auch in anderen Klassen als sehr n\u00fctzlich heraus um Tasks sp\u00e4ter auszuf\u00fchren. '\n$run3 = 'Videos k\u00f6nnen leider immernoch nicht abgespielt werden.'\n\n$d = $word.ActiveDocument\n\n# Find the paragraph whose whole text is exactly \"Tutorial\" (the short\n# diary placeholder written on \"6. Mai\"). Remember its start/end so we can\n# find it again unambiguously later (the new text we insert also contains\n# the word \"Tutorial\", so a second text search would be ambiguous).\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Tutorial\")\nif (-not $found) {\n    throw 'Could not find \"Tutorial\" text to replace.'\n}\n$oldStart = $rng.Start\n$oldEnd = $rng.End\n\n# Build a tiny OOXML fragment with three independent <w:r> runs (no run\n# properties, same as the original single run) and insert it as a new\n# range collapsed to the start of \"Tutorial\" - this puts the new runs\n# right before the old \"Tutorial\" run without disturbing the paragraph's\n# own formatting (<w:pPr>) or the bookmark that follows the old run.\n$insertionPoint = $d.Range($oldStart, $oldStart)\n$ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>' + $run1 + '</w:t></w:r><w:r><w:t xml:space=\"preserve\">' + $run2 + '</w:t></w:r><w:r><w:t>' + $run3 + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertionPoint.InsertXML($ooxml) | Out-Null\n\n# Now clear out the original \"Tutorial\" text (shifted right by the length\n# of the text we just inserted in front of it), leaving the bookmark\n# (_GoBack) that immediately follows it untouched and in place.\n$insertedLength = $run1.Length + $run2.Length + $run3.Length\n$oldRng = $d.Range($oldStart + $insertedLength, $oldEnd + $insertedLength)\n$oldRng.Text = \"\"\n"}
